$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.771.37"
$ws.Range("E2").Value = "  +5.84%  "
$ws.Range("D3").Value = "3.530.53"
$ws.Range("E3").Value = "  +9.10%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "568.88"
$ws.Range("E5").Value = "  +7.52%  "
$ws.Range("D6").Value = "187.97"
$ws.Range("E6").Value = "  +9.85%  "
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  +4.36%  "
$ws.Range("D8").Value = "3.523.52"
$ws.Range("E8").Value = "  +9.03%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "0.633"
$ws.Range("E10").Value = "  +4.75%  "
$ws.Range("E11").Value = "  +13.12%  "
$ws.Range("D12").Value = "54.70"
$ws.Range("E12").Value = "  +3.18%  "
$ws.Range("E13").Value = "  +6.43%  "
$ws.Range("D14").Value = "9.40"
$ws.Range("E14").Value = "  +3.15%  "
$ws.Range("D15").Value = "4.100.03"
$ws.Range("E15").Value = "  +9.20%  "
$ws.Range("D16").Value = "3.536.27"
$ws.Range("E16").Value = "  +8.99%  "
$ws.Range("E17").Value = "  +4.73%  "
$ws.Range("D18").Value = "66.838.69"
$ws.Range("D19").Value = "18.24"
$ws.Range("E19").Value = "  +6.10%  "
$ws.Range("D20").Value = "12.00"
$ws.Range("E20").Value = "  +8.60%  "
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").Value = "  +3.11%  "
$ws.Range("D22").Value = "427.32"
$ws.Range("E22").Value = "  +16.75%  "
$ws.Range("D23").Value = "4.20"
$ws.Range("E23").Value = "  +11.87%  "
$ws.Range("D24").Value = "85.10"
$ws.Range("E24").Value = "  +5.16%  "
$ws.Range("D25").Value = "4.12"
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("D26").Value = "11.15"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "2.89"
$ws.Range("E27").Value = "  +9.76%  "
$ws.Range("D28").Value = "12.26"
$ws.Range("E28").Value = "  +9.08%  "
$ws.Range("D29").Value = "9.24"
$ws.Range("E29").Value = "  +12.96%  "
$ws.Range("D30").Value = "30.29"
$ws.Range("E30").Value = "  +6.54%  "
$ws.Range("D31").Value = "643.45"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").Value = "6.61"
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("D33").Value = "11.72"
$ws.Range("E33").Value = "  +4.88%  "
$ws.Range("D35").Value = "59.69"
$ws.Range("E35").Value = "  +5.56%  "
$ws.Range("D36").Value = "38.52"
$ws.Range("E36").Value = "  +5.45%  "
$ws.Range("D37").Value = "0.148"
$ws.Range("E37").Value = "  +20.55%  "
$ws.Range("D38").Value = "0.0₃0812"
$ws.Range("E38").Value = "  +13.43%  "
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").Value = "0.392"
$ws.Range("E40").Value = "  +3.97%  "
$ws.Range("E41").Value = "  +14.38%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "3.044.99"
$ws.Range("E43").Value = "  +5.74%  "
$ws.Range("E44").Value = "  +4.87%  "
$ws.Range("D45").Value = "2.88"
$ws.Range("E45").Value = "  +11.51%  "
$ws.Range("D46").Value = "3.32"
$ws.Range("E46").Value = "  +7.71%  "
$ws.Range("D47").Value = "0.0420"
$ws.Range("E47").Value = "  +6.77%  "
$ws.Range("D48").Value = "2.76"
$ws.Range("E48").Value = "  +2.77%  "
$ws.Range("E49").Value = "  +5.58%  "
$ws.Range("D50").Value = "141.51"
$ws.Range("E50").Value = "  +5.70%  "
$ws.Range("D51").Value = "8.63"
$ws.Range("E51").Value = "  +10.71%  "
